# The sheet already has one "~TFM_INS / FLO_SHAR" row for year 2030
# (row 5). This change extends the same flow-share setting out to the
# remaining model years (2035, 2040, 2045, 2050) by duplicating that
# row's values/formatting and only changing the Year column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B5:H5").Copy($ws.Range("B6:H6"))
$ws.Range("B5:H5").Copy($ws.Range("B7:H7"))
$ws.Range("B5:H5").Copy($ws.Range("B8:H8"))
$ws.Range("B5:H5").Copy($ws.Range("B9:H9"))

$ws.Range("F6").Value = 2035
$ws.Range("F7").Value = 2040
$ws.Range("F8").Value = 2045
$ws.Range("F9").Value = 2050

# Leave the selection where the author's saved file shows it.
$ws.Range("F10").Select()
